# Mart Management System - slide.pptx
# "chore: assign member in slide"
#
# Re-assigns the presenter/member name shown on several "member" slides.
# Each target shape is a simple one-paragraph / one-run text box whose
# a:bodyPr uses <a:spAutoFit/>; this runtime recomputes the shape's
# Height as soon as TextRange.Text is written (even though real
# PowerPoint did not need to resize these particular boxes for this
# edit), so we restore the original Height immediately afterwards to
# keep the shape geometry byte-identical to the source file.

$p = $ppt.ActivePresentation

# EMU -> Points helper tuned for this runtime's float32 truncation of
# Shape.Height/Width (pt is cast to float32 then multiplied by 12700 and
# truncated). Adding 0.5 EMU before converting lands safely in the middle
# of the rounding bucket so the value always truncates back to $emu.
function PtForEmu([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

$origHeightPt = PtForEmu 400200

function Set-MemberName($slideIndex, $shapeIndex, $text) {
    $shp = $p.Slides.Item($slideIndex).Shapes.Item($shapeIndex)
    $shp.TextFrame.TextRange.Text = $text
    $shp.Height = $origHeightPt
}

Set-MemberName 14 2 "លិ ស្រីម៉ា​ ​- Le Sreyma"
Set-MemberName 4  7 "ធិន​ សុីវធាន - Then Sivthean"
Set-MemberName 5  5 "ធិន​ សុីវធាន - Then Sivthean"
Set-MemberName 6  3 "អន ភក្តី - Orn Pheakdey"
Set-MemberName 7  2 "អន ភក្តី - Orn Pheakdey"
Set-MemberName 8  2 "អន ភក្តី - Orn Pheakdey"
Set-MemberName 9  2 "អន ភក្តី - Orn Pheakdey"
